$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection / active cell
$ws.Range("F22").Select() | Out-Null

# Update values in "test Find" block (rows 3-7)
$ws.Range("C3").Value = 16.848
$ws.Range("D3").Value = 159.078
$ws.Range("C4").Value = 18.535
$ws.Range("D4").Value = 186.143
$ws.Range("C5").Value = 25.616
$ws.Range("D5").Value = 247.488
$ws.Range("C6").Value = 43.214
$ws.Range("D6").Value = 326.721
$ws.Range("C7").Value = 71.015
$ws.Range("D7").Value = 376.969

# Update values in "test Insert" block (rows 11-15)
$ws.Range("C11").Value = 126.142
$ws.Range("D11").Value = 160.545
$ws.Range("C12").Value = 165.859
$ws.Range("D12").Value = 198.403
$ws.Range("C13").Value = 264.426
$ws.Range("D13").Value = 273.973
$ws.Range("C14").Value = 370.57
$ws.Range("D14").Value = 345.289
$ws.Range("C15").Value = 439.923
$ws.Range("D15").Value = 419.653

# Update values in "test Remove" block (rows 19-23)
$ws.Range("C19").Value = 557.982
$ws.Range("D19").Value = 646.011
$ws.Range("C20").Value = 576.094
$ws.Range("D20").Value = 692.218
$ws.Range("C21").Value = 657.282
$ws.Range("D21").Value = 783.454
$ws.Range("C22").Value = 790.432
$ws.Range("D22").Value = 917.259
$ws.Range("C23").Value = 872.337
$ws.Range("D23").Value = 1011.432
